$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Summary sheet - update headline stats after trade #61 (sheet row 90/20)
#    closed and trade #118 (new) opened.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.13
$summary.Range("B4").Value = 0.24
$summary.Range("B5").Value = 0.05
$summary.Range("B6").Value = 89
$summary.Range("B8").Value = 36
$summary.Range("B9").Value = 49.44

# ---------------------------------------------------------------------------
# 2) Strategy Status sheet - update the "momentum" strategy row (row 11)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C11").Value = 99.2
$status.Range("D11").Value = 19
$status.Range("E11").Value = -0.8
$status.Range("F11").Value = -0.8
$status.Range("G11").Value = 21.05

# ---------------------------------------------------------------------------
# 3) All Trades sheet - close trade #89 (row 90) and append new trade #118
#    (row 119)
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Range("G90").Value = 0.96
$allTrades.Range("H90").Value = "CLOSED"
$allTrades.Range("I90").Value = -1.0309
$allTrades.Range("J90").Value = -0.01
$allTrades.Range("K90").Value = 99.2
$allTrades.Range("L90").Value = "early_exit"
$allTrades.Range("M90").Value = 0.13

$allTrades.Range("A119").Value = 118
$allTrades.Range("B119").Value = "'2026-02-18"
$allTrades.Range("C119").Value = "00:23:13"
$allTrades.Range("D119").Value = "MarketMaking"
$allTrades.Range("E119").Value = "DOWN"
$allTrades.Range("F119").Value = 0.97
$allTrades.Range("H119").Value = "OPEN"
$allTrades.Range("I119").Value = 0
$allTrades.Range("J119").Value = 0
$allTrades.Range("K119").Value = 99.410254715139
$allTrades.Range("M119").Value = 0
$allTrades.Range("N119").Value = 0
$allTrades.Range("O119").Value = 0
$allTrades.Range("P119").Value = 0.6
$allTrades.Range("Q119").Value = "Normal spread capture: 198 bps"

# ---------------------------------------------------------------------------
# 4) momentum sheet - close trade #89 (row 20)
# ---------------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Range("G20").Value = 0.96
$momentum.Range("H20").Value = "CLOSED"
$momentum.Range("I20").Value = -1.0309
$momentum.Range("J20").Value = -0.01
$momentum.Range("K20").Value = 99.2
$momentum.Range("P20").Value = "early_exit"
$momentum.Range("Q20").Value = 0.13

# ---------------------------------------------------------------------------
# 5) MarketMaking sheet - append new trade #118 (row 39)
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Range("A39").Value = 118
$marketMaking.Range("B39").Value = "'2026-02-18"
$marketMaking.Range("C39").Value = "00:23:13"
$marketMaking.Range("D39").Value = "MarketMaking"
$marketMaking.Range("E39").Value = "DOWN"
$marketMaking.Range("F39").Value = 0.97
$marketMaking.Range("H39").Value = "OPEN"
$marketMaking.Range("I39").Value = 0
$marketMaking.Range("J39").Value = 0
$marketMaking.Range("K39").Value = 99.410254715139
$marketMaking.Range("L39").Value = 0
$marketMaking.Range("M39").Value = 0
$marketMaking.Range("N39").Value = 0.6
$marketMaking.Range("O39").Value = "Normal spread capture: 198 bps"
$marketMaking.Range("Q39").Value = 0
